$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '22.393.77'
$ws.Cells.Item(2, 5).Value = '  -0.15%  '

$ws.Cells.Item(3, 4).Value = '1.565.35'
$ws.Cells.Item(3, 5).Value = '  -0.07%  '

$ws.Cells.Item(5, 5).Value = '  -0.14%  '

$ws.Cells.Item(6, 4).Value = '286.05'
$ws.Cells.Item(6, 5).Value = '  +0.37%  '

$ws.Cells.Item(7, 4).Value = '0.3719'
$ws.Cells.Item(7, 5).Value = '  +2.47%  '

$ws.Cells.Item(8, 2).Value = 'OKB'
$ws.Cells.Item(8, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(8, 4).Value = '46.49'
$ws.Cells.Item(8, 5).Value = '  -4.01%  '

$ws.Cells.Item(9, 2).Value = 'Cardano'
$ws.Cells.Item(9, 3).Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Cells.Item(9, 4).Value = '0.3272'
$ws.Cells.Item(9, 5).Value = '  -1.41%  '

$ws.Cells.Item(10, 4).Value = '1.143'
$ws.Cells.Item(10, 5).Value = '  +1.99%  '

$ws.Cells.Item(11, 4).Value = '0.07406'
$ws.Cells.Item(11, 5).Value = '  +0.33%  '

$ws.Cells.Item(12, 5).Value = '  -0.08%  '

$ws.Cells.Item(13, 4).Value = '20.35'
$ws.Cells.Item(13, 5).Value = '  -1.94%  '

$ws.Cells.Item(14, 4).Value = '5.833'
$ws.Cells.Item(14, 5).Value = '  -1.79%  '

$ws.Cells.Item(15, 4).Value = '6.817'
$ws.Cells.Item(15, 5).Value = '  -1.04%  '

$ws.Cells.Item(16, 4).Value = '1.564.32'
$ws.Cells.Item(16, 5).Value = '  -0.13%  '

$ws.Cells.Item(17, 4).Value = '0.00001095'
$ws.Cells.Item(17, 5).Value = '  -0.71%  '

$ws.Cells.Item(18, 4).Value = '0.06690'
$ws.Cells.Item(18, 5).Value = '  -0.34%  '

$ws.Cells.Item(19, 4).Value = '86.03'
$ws.Cells.Item(19, 5).Value = '  -2.05%  '

$ws.Cells.Item(20, 4).Value = '1.0000'
$ws.Cells.Item(20, 5).Value = '  -0.20%  '

$ws.Cells.Item(21, 4).Value = '6.325'
$ws.Cells.Item(21, 5).Value = '  -0.09%  '

$ws.Cells.Item(22, 5).Value = '  +0.30%  '

$ws.Cells.Item(23, 4).Value = '11.75'
$ws.Cells.Item(23, 5).Value = '  -1.92%  '

$ws.Cells.Item(24, 4).Value = '22.387.25'
$ws.Cells.Item(24, 5).Value = '  -0.16%  '

$ws.Cells.Item(25, 4).Value = '2.283'
$ws.Cells.Item(25, 5).Value = '  -3.92%  '

$ws.Cells.Item(26, 4).Value = '2.558'
$ws.Cells.Item(26, 5).Value = '  +1.00%  '

$ws.Cells.Item(27, 4).Value = '151.10'
$ws.Cells.Item(27, 5).Value = '  +0.46%  '

$ws.Cells.Item(28, 4).Value = '19.27'
$ws.Cells.Item(28, 5).Value = '  -0.61%  '

$ws.Cells.Item(29, 4).Value = '4.937'
$ws.Cells.Item(29, 5).Value = '  -1.25%  '

$ws.Cells.Item(30, 4).Value = '123.65'
$ws.Cells.Item(30, 5).Value = '  -0.02%  '

$ws.Cells.Item(31, 4).Value = '1.741.03'
$ws.Cells.Item(31, 5).Value = '  -0.13%  '

$ws.Cells.Item(32, 4).Value = '1.050'
$ws.Cells.Item(32, 5).Value = '  +1.88%  '

$ws.Cells.Item(33, 4).Value = '1.952'
$ws.Cells.Item(33, 5).Value = '  -2.71%  '

$ws.Cells.Item(34, 4).Value = '5.932'
$ws.Cells.Item(34, 5).Value = '  -2.50%  '

$ws.Cells.Item(35, 4).Value = '9.618'
$ws.Cells.Item(35, 5).Value = '  -0.97%  '

$ws.Cells.Item(36, 4).Value = '0.08204'
$ws.Cells.Item(36, 5).Value = '  -0.41%  '

$ws.Cells.Item(37, 4).Value = '1.322'
$ws.Cells.Item(37, 5).Value = '  +2.74%  '

$ws.Cells.Item(38, 4).Value = '0.02372'
$ws.Cells.Item(38, 5).Value = '  -1.43%  '

$ws.Cells.Item(39, 4).Value = '0.06283'
$ws.Cells.Item(39, 5).Value = '  -2.45%  '

$ws.Cells.Item(40, 4).Value = '0.2181'
$ws.Cells.Item(40, 5).Value = '  -2.43%  '

$ws.Cells.Item(41, 4).Value = '5.222'
$ws.Cells.Item(41, 5).Value = '  -2.56%  '

$ws.Cells.Item(42, 4).Value = '11.09'
$ws.Cells.Item(42, 5).Value = '  -0.73%  '

$ws.Cells.Item(43, 4).Value = '0.6098'
$ws.Cells.Item(43, 5).Value = '  -2.69%  '

$ws.Cells.Item(44, 4).Value = '1.000'
$ws.Cells.Item(44, 5).Value = '  -0.09%  '

$ws.Cells.Item(45, 4).Value = '13.74'
$ws.Cells.Item(45, 5).Value = '  +0.10%  '

$ws.Cells.Item(46, 5).Value = '  -1.91%  '

$ws.Cells.Item(47, 4).Value = '3.742'
$ws.Cells.Item(47, 5).Value = '  -0.21%  '

$ws.Cells.Item(48, 4).Value = '2.002'
$ws.Cells.Item(48, 5).Value = '  -1.13%  '

$ws.Cells.Item(49, 4).Value = '123.69'
$ws.Cells.Item(49, 5).Value = '  +0.35%  '

$ws.Cells.Item(50, 4).Value = '1.176'
$ws.Cells.Item(50, 5).Value = '  -2.83%  '

$ws.Cells.Item(51, 5).Value = '  -0.68%  '
